$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H113").Value = 1672.7667
$ws.Range("I113").Value = 1244.2727
$ws.Range("J113").Value = 1920.8422
$ws.Range("K113").Value = 1244.2727
$ws.Range("L113").Value = 1920.8422
$ws.Range("M113").Value = 2009.7273
$ws.Range("N113").Value = -8428.842199999999

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 2138942.5
$ws.Range("I116").Value = 9617159
$ws.Range("J116").Value = 2309.1428
$ws.Range("K116").Value = 9617159
$ws.Range("L116").Value = 2309.1428
$ws.Range("M116").Value = -9613717
$ws.Range("N116").Value = -9193.1428

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 685
$ws.Range("I129").Value = 264.42856
$ws.Range("K129").Value = 793.28568
$ws.Range("M129").Value = 4206.71432

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 2473.9592
$ws.Range("I132").Value = 2409.2888
$ws.Range("K132").Value = 7227.866399999999
$ws.Range("M132").Value = -4697.866399999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 14506693
$ws.Range("I32").Value = 20837846
$ws.Range("J32").Value = 35485.668
$ws.Range("K32").Value = 20837846
$ws.Range("L32").Value = 35485.668
$ws.Range("M32").Value = -20837559
$ws.Range("N32").Value = -36059.668

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1422.2142
$ws.Range("I45").Value = 973.5833
$ws.Range("J45").Value = 4114
$ws.Range("K45").Value = 973.5833
$ws.Range("L45").Value = 4114
$ws.Range("M45").Value = -596.5833
$ws.Range("N45").Value = -4868

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 669.1667
$ws.Range("I97").Value = 347.25
$ws.Range("J97").Value = 830.125
$ws.Range("K97").Value = 347.25
$ws.Range("L97").Value = 830.125
$ws.Range("M97").Value = 148.75
$ws.Range("N97").Value = -1822.125

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 1441.1562
$ws.Range("I132").Value = 900.4231
$ws.Range("J132").Value = 3784.3333
$ws.Range("K132").Value = 2701.2693
$ws.Range("L132").Value = 11352.9999
$ws.Range("M132").Value = -171.2692999999999
$ws.Range("N132").Value = -16412.9999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 2899.6316
$ws.Range("I105").Value = 1205
$ws.Range("J105").Value = 2993.7778
$ws.Range("K105").Value = 1205
$ws.Range("L105").Value = 2993.7778
$ws.Range("M105").Value = 542
$ws.Range("N105").Value = -6487.7778

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 1412.9756
$ws.Range("I134").Value = 1124.8788
$ws.Range("J134").Value = 2601.375
$ws.Range("K134").Value = 3374.6364
$ws.Range("L134").Value = 7804.125
$ws.Range("M134").Value = -839.6363999999999
$ws.Range("N134").Value = -12874.125

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4664.278
$ws.Range("I31").Value = 3367.5454
$ws.Range("J31").Value = 6702
$ws.Range("K31").Value = 3367.5454
$ws.Range("L31").Value = 6702
$ws.Range("M31").Value = -3072.5454
$ws.Range("N31").Value = -7292

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 4664.278
$ws.Range("I34").Value = 3367.5454
$ws.Range("J34").Value = 6702
$ws.Range("K34").Value = 3367.5454
$ws.Range("L34").Value = 6702
$ws.Range("M34").Value = -3165.5454
$ws.Range("N34").Value = -7106

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2817.348
$ws.Range("I58").Value = 1440.7693
$ws.Range("K58").Value = 1440.7693
$ws.Range("M58").Value = -1237.7693

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 2430.9
$ws.Range("I99").Value = 1756.4
$ws.Range("J99").Value = 3105.4
$ws.Range("K99").Value = 1756.4
$ws.Range("L99").Value = 3105.4
$ws.Range("M99").Value = -258.4000000000001
$ws.Range("N99").Value = -6101.4

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 2430.9
$ws.Range("I126").Value = 1756.4
$ws.Range("J126").Value = 3105.4
$ws.Range("K126").Value = 5269.200000000001
$ws.Range("L126").Value = 9316.200000000001
$ws.Range("M126").Value = -2799.200000000001
$ws.Range("N126").Value = -14256.2

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 2817.348
$ws.Range("I136").Value = 1440.7693
$ws.Range("K136").Value = 4322.3079
$ws.Range("M136").Value = -1772.3079

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H121").Value = 12501144
$ws.Range("I121").Value = 300
$ws.Range("J121").Value = 16668092
$ws.Range("K121").Value = 900
$ws.Range("L121").Value = 50004276
$ws.Range("M121").Value = 410
$ws.Range("N121").Value = -50006896

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1569.0714
$ws.Range("I102").Value = 1492.4546
$ws.Range("J102").Value = 1850
$ws.Range("K102").Value = 1492.4546
$ws.Range("L102").Value = 1850
$ws.Range("M102").Value = 129.5454
$ws.Range("N102").Value = -5094

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 14355.733
$ws.Range("I113").Value = 1844.4
$ws.Range("J113").Value = 20611.4
$ws.Range("K113").Value = 1844.4
$ws.Range("L113").Value = 20611.4
$ws.Range("M113").Value = 325.5999999999999
$ws.Range("N113").Value = -24951.4

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 4764011.5
$ws.Range("I126").Value = 8334701
$ws.Range("K126").Value = 25004103
$ws.Range("M126").Value = -25001633

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2664.2144
$ws.Range("I7").Value = 2499.8572
$ws.Range("J7").Value = 2828.5715
$ws.Range("K7").Value = 2499.8572
$ws.Range("L7").Value = 2828.5715
$ws.Range("M7").Value = -2387.8572
$ws.Range("N7").Value = -3052.5715

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 2343.7144
$ws.Range("I82").Value = 1400.4
$ws.Range("K82").Value = 1400.4
$ws.Range("M82").Value = -1039.4

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H85").Value = 2343.7144
$ws.Range("I85").Value = 1400.4
$ws.Range("K85").Value = 1400.4
$ws.Range("M85").Value = -152.4000000000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 3543.0862
$ws.Range("I122").Value = 3524.5405
$ws.Range("K122").Value = 10573.6215
$ws.Range("M122").Value = -8123.621500000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 2664.2144
$ws.Range("I126").Value = 2499.8572
$ws.Range("J126").Value = 2828.5715
$ws.Range("K126").Value = 7499.571599999999
$ws.Range("L126").Value = 8485.7145
$ws.Range("M126").Value = -5029.571599999999
$ws.Range("N126").Value = -13425.7145

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1407.1428
$ws.Range("I81").Value = 1075
$ws.Range("J81").Value = 1540
$ws.Range("K81").Value = 2150
$ws.Range("L81").Value = 3080
$ws.Range("M81").Value = -1089
$ws.Range("N81").Value = -5202

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H84").Value = 1407.1428
$ws.Range("I84").Value = 1075
$ws.Range("J84").Value = 1540
$ws.Range("K84").Value = 10750
$ws.Range("L84").Value = 15400
$ws.Range("M84").Value = -5446
$ws.Range("N84").Value = -26008

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 2348.88
$ws.Range("I126").Value = 1668.875
$ws.Range("J126").Value = 3557.7778
$ws.Range("K126").Value = 5006.625
$ws.Range("L126").Value = 10673.3334
$ws.Range("M126").Value = -2536.625
$ws.Range("N126").Value = -15613.3334

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 4002.389
$ws.Range("J136").Value = 2507
$ws.Range("L136").Value = 7521
$ws.Range("N136").Value = -12621
